$wb = $excel.ActiveWorkbook

# --- Sheet: ELF-bldg-winter ---
$wsWinter = $wb.Worksheets.Item("ELF-bldg-winter")
$wsWinter.Range("D2").Value = 1.35615
$wsWinter.Range("D5").Value = 1.40111
$wsWinter.Range("D7").Value = 1.40111

# --- Sheet: ELF-bldg-summer ---
$wsSummer = $wb.Worksheets.Item("ELF-bldg-summer")
$wsSummer.Range("B3").Value = 3.15878
$wsSummer.Range("D3").Value = 2.101

# --- Sheet: ELF-vehicles ---
$wsVehicles = $wb.Worksheets.Item("ELF-vehicles")
$wsVehicles.Range("B4").Value = 1.10177
$wsVehicles.Range("C4").Value = 1.03216
$wsVehicles.Range("B5").Value = 1.10177
$wsVehicles.Range("C5").Value = 1.03216
$wsVehicles.Range("B6").Value = 1.10177
$wsVehicles.Range("C6").Value = 1.03216
$wsVehicles.Range("B7").Value = 1.10177
$wsVehicles.Range("C7").Value = 1.03216
